# Update "想去人数" (F column) figures across the four worksheets to the
# newly published numbers (regenerated gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 4441
$ws.Range("F7").Value = 3547
$ws.Range("F8").Value = 1011
$ws.Range("F11").Value = 323
$ws.Range("F12").Value = 321
$ws.Range("F13").Value = 2411
$ws.Range("F15").Value = 32
$ws.Range("F18").Value = 534
$ws.Range("F19").Value = 252
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 10018
$ws.Range("F22").Value = 5952
$ws.Range("F23").Value = 382
$ws.Range("F25").Value = 824
$ws.Range("F31").Value = 459
$ws.Range("F32").Value = 112
$ws.Range("F33").Value = 240
$ws.Range("F36").Value = 4816
$ws.Range("F40").Value = 25
$ws.Range("F41").Value = 64

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 22
$ws.Range("F15").Value = 3528
$ws.Range("F18").Value = 7

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8705
$ws.Range("F3").Value = 419
$ws.Range("F4").Value = 1562

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8705
$ws.Range("F4").Value = 419
$ws.Range("F5").Value = 1562
$ws.Range("F7").Value = 4441
$ws.Range("F10").Value = 3547
$ws.Range("F11").Value = 1011
$ws.Range("F13").Value = 321
$ws.Range("F14").Value = 2411
$ws.Range("F17").Value = 22
$ws.Range("F20").Value = 32
$ws.Range("F23").Value = 534
$ws.Range("F24").Value = 252
$ws.Range("F25").Value = 10018
$ws.Range("F26").Value = 3528
$ws.Range("F28").Value = 382
$ws.Range("F30").Value = 824
$ws.Range("F35").Value = 459
$ws.Range("F36").Value = 112
$ws.Range("F37").Value = 240
$ws.Range("F39").Value = 7
$ws.Range("F40").Value = 4816
$ws.Range("F45").Value = 64

$wb.Save()
